$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.90017149061062
$ws.Range("C2").Value = 7.450482347439462
$ws.Range("D2").Value = 13.36095203485508
$ws.Range("E2").Value = 13.60012849097872
$ws.Range("G2").Value = 3.712581000672822
$ws.Range("J2").Value = 8.355446309843638
$ws.Range("L2").Value = 12.30066027045948
$ws.Range("M2").Value = 19.10344080315191
$ws.Range("N2").Value = 20.64716776482446
$ws.Range("O2").Value = 32.71741759991863

$ws.Range("B3").Value = 20.61256484019727
$ws.Range("C3").Value = 7.188280921603954
$ws.Range("D3").Value = 13.37682016901493
$ws.Range("E3").Value = 13.63232731925872
$ws.Range("G3").Value = 3.715181920250302
$ws.Range("J3").Value = 8.352952483653093
$ws.Range("L3").Value = 12.30938497058558
$ws.Range("M3").Value = 19.05135426537491
$ws.Range("N3").Value = 20.71579964130773
$ws.Range("O3").Value = 32.78131793631842

$ws.Range("B4").Value = 20.43810243344332
$ws.Range("C4").Value = 7.020948454680127
$ws.Range("D4").Value = 13.38875353698746
$ws.Range("E4").Value = 13.65332109316314
$ws.Range("G4").Value = 3.716864284877447
$ws.Range("J4").Value = 8.351439606168427
$ws.Range("L4").Value = 12.31626747406728
$ws.Range("M4").Value = 19.02228696092054
$ws.Range("N4").Value = 20.75988118260295
$ws.Range("O4").Value = 32.82745567913754

$ws.Range("B5").Value = 20.36762754702159
$ws.Range("C5").Value = 6.951227051167436
$ws.Range("D5").Value = 13.39416727251472
$ws.Range("E5").Value = 13.66218454573409
$ws.Range("G5").Value = 3.717571405014507
$ws.Range("J5").Value = 8.350827744354616
$ws.Range("L5").Value = 12.31945638928321
$ws.Range("M5").Value = 19.01118276326843
$ws.Range("N5").Value = 20.77833450825413
$ws.Range("O5").Value = 32.84798929791059

$ws.Range("B6").Value = 20.35596514454025
$ws.Range("C6").Value = 6.939559286745091
$ws.Range("D6").Value = 13.39509948427631
$ws.Range("E6").Value = 13.66367495602004
$ws.Range("G6").Value = 3.717690125000233
$ws.Range("J6").Value = 8.350726426940755
$ws.Range("L6").Value = 12.32000913073662
$ws.Range("M6").Value = 19.00938389964197
$ws.Range("N6").Value = 20.78142829464168
$ws.Range("O6").Value = 32.85150339210311

$ws.Range("B7").Value = 20.43714936320672
$ws.Range("C7").Value = 7.020014283442737
$ws.Range("D7").Value = 13.38882431842478
$ws.Range("E7").Value = 13.65343937955078
$ws.Range("G7").Value = 3.716873734026272
$ws.Range("J7").Value = 8.351431335518456
$ws.Range("L7").Value = 12.31630892432796
$ws.Range("M7").Value = 19.02213419533103
$ws.Range("N7").Value = 20.76012806558119
$ws.Range("O7").Value = 32.82772559454513

$ws.Range("B8").Value = 20.80061371811502
$ws.Range("C8").Value = 7.361426206589031
$ws.Range("D8").Value = 13.36596882596353
$ws.Range("E8").Value = 13.61097708418091
$ws.Range("G8").Value = 3.713460117959486
$ws.Range("J8").Value = 8.354582652559451
$ws.Range("L8").Value = 12.30335232963927
$ws.Range("M8").Value = 19.08488166024382
$ws.Range("N8").Value = 20.67043007944358
$ws.Range("O8").Value = 32.73801571210094

$ws.Range("B9").Value = 21.52652581729663
$ws.Range("C9").Value = 7.978250897534343
$ws.Range("D9").Value = 13.33852360056967
$ws.Range("E9").Value = 13.5373879114069
$ws.Range("G9").Value = 3.707440318630014
$ws.Range("J9").Value = 8.360906646773428
$ws.Range("L9").Value = 12.29001911481096
$ws.Range("M9").Value = 19.23066502699053
$ws.Range("N9").Value = 20.50986243843492
$ws.Range("O9").Value = 32.61701344819159

$ws.Range("B10").Value = 22.0628915007718
$ws.Range("C10").Value = 8.396577621804688
$ws.Range("D10").Value = 13.32894061333142
$ws.Range("E10").Value = 13.4891822038106
$ws.Range("G10").Value = 3.703424114687132
$ws.Range("J10").Value = 8.365639598929901
$ws.Range("L10").Value = 12.28754231264595
$ws.Range("M10").Value = 19.35107968468857
$ws.Range("N10").Value = 20.40113454812478
$ws.Range("O10").Value = 32.56177688497115

$ws.Range("B11").Value = 22.30652086646069
$ws.Range("C11").Value = 8.5788536912838
$ws.Range("D11").Value = 13.3268738771464
$ws.Range("E11").Value = 13.46851615685673
$ws.Range("G11").Value = 3.701684353933751
$ws.Range("J11").Value = 8.367811325062098
$ws.Range("L11").Value = 12.28799454036597
$ws.Range("M11").Value = 19.40861945009618
$ws.Range("N11").Value = 20.35365576827843
$ws.Range("O11").Value = 32.54398916753362

$ws.Range("B12").Value = 22.3986401676916
$ws.Range("C12").Value = 8.646691145130083
$ws.Range("D12").Value = 13.3264202747207
$ws.Range("E12").Value = 13.46087144046312
$ws.Range("G12").Value = 3.701038022197114
$ws.Range("J12").Value = 8.368636354742851
$ws.Range("L12").Value = 12.28839181718043
$ws.Range("M12").Value = 19.43079365424249
$ws.Range("N12").Value = 20.33596015535534
$ws.Range("O12").Value = 32.53831059610829

$ws.Range("B13").Value = 22.37880808873618
$ws.Range("C13").Value = 8.632134383033376
$ws.Range("D13").Value = 13.32650334442145
$ws.Range("E13").Value = 13.46250982253612
$ws.Range("G13").Value = 3.701176667384338
$ws.Range("J13").Value = 8.368458552922874
$ws.Range("L13").Value = 12.28829622094254
$ws.Range("M13").Value = 19.42600111130431
$ws.Range("N13").Value = 20.33975863329847
$ws.Range("O13").Value = 32.53948653269431

$ws.Range("B14").Value = 22.31410281868369
$ws.Range("C14").Value = 8.5844586717279
$ws.Range("D14").Value = 13.32682996904507
$ws.Range("E14").Value = 13.46788359568745
$ws.Range("G14").Value = 3.701630930131881
$ws.Range("J14").Value = 8.367879148003331
$ws.Range("L14").Value = 12.28802270075636
$ws.Range("M14").Value = 19.41043608747498
$ws.Range("N14").Value = 20.35219426464268
$ws.Range("O14").Value = 32.54350079045474

$ws.Range("B15").Value = 22.27444853387766
$ws.Range("C15").Value = 8.555100501182155
$ws.Range("D15").Value = 13.32707286329657
$ws.Range("E15").Value = 13.47119875109473
$ws.Range("G15").Value = 3.70191080223006
$ws.Range("J15").Value = 8.367524588834966
$ws.Range("L15").Value = 12.28788456562704
$ws.Range("M15").Value = 19.40095185021395
$ws.Range("N15").Value = 20.3598483343291
$ws.Range("O15").Value = 32.54609737029809

$ws.Range("B16").Value = 22.04695504238467
$ws.Range("C16").Value = 8.384501098644765
$ws.Range("D16").Value = 13.32912176028365
$ws.Range("E16").Value = 13.49055814588446
$ws.Range("G16").Value = 3.703539561815725
$ws.Range("J16").Value = 8.365498043849188
$ws.Range("L16").Value = 12.28754444522701
$ws.Range("M16").Value = 19.34737389172406
$ws.Range("N16").Value = 20.40427715834516
$ws.Range("O16").Value = 32.56308721309189

$ws.Range("B17").Value = 21.9072400032186
$ws.Range("C17").Value = 8.277764204025367
$ws.Range("D17").Value = 13.33096546713157
$ws.Range("E17").Value = 13.50275759334538
$ws.Range("G17").Value = 3.704561047974116
$ws.Range("J17").Value = 8.364259579200144
$ws.Range("L17").Value = 12.28773949409277
$ws.Range("M17").Value = 19.31520465021313
$ws.Range("N17").Value = 20.432039425851
$ws.Range("O17").Value = 32.5753912390358

$ws.Range("B18").Value = 21.82685010794942
$ws.Range("C18").Value = 8.215618474614235
$ws.Range("D18").Value = 13.3322417576568
$ws.Range("E18").Value = 13.50989330188196
$ws.Range("G18").Value = 3.70515679424673
$ws.Range("J18").Value = 8.363549008148421
$ws.Range("L18").Value = 12.28800036409637
$ws.Range("M18").Value = 19.29696253371316
$ws.Range("N18").Value = 20.44819419671193
$ws.Range("O18").Value = 32.58315895495688

$ws.Range("B19").Value = 21.79962904017308
$ws.Range("C19").Value = 8.194448661915398
$ws.Range("D19").Value = 13.33271097579712
$ws.Range("E19").Value = 13.51232976970734
$ws.Range("G19").Value = 3.70535991655012
$ws.Range("J19").Value = 8.363308724743488
$ws.Range("L19").Value = 12.28811425650671
$ws.Range("M19").Value = 19.2908312163054
$ws.Range("N19").Value = 20.4536960260968
$ws.Range("O19").Value = 32.58590754676967

$ws.Range("B20").Value = 21.92211658585304
$ws.Range("C20").Value = 8.289204764537974
$ws.Range("D20").Value = 13.33074686712201
$ws.Range("E20").Value = 13.50144663868456
$ws.Range("G20").Value = 3.70445145934054
$ws.Range("J20").Value = 8.36439123391953
$ws.Range("L20").Value = 12.28770334982075
$ws.Range("M20").Value = 19.31860221716798
$ws.Range("N20").Value = 20.42906477814737
$ws.Range("O20").Value = 32.57400994811731

$ws.Range("B21").Value = 22.33311273281668
$ws.Range("C21").Value = 8.598494625788016
$ws.Range("D21").Value = 13.32672510767155
$ws.Range("E21").Value = 13.46630027846044
$ws.Range("G21").Value = 3.701497164026289
$ws.Range("J21").Value = 8.368049261853214
$ws.Range("L21").Value = 12.28809691445444
$ws.Range("M21").Value = 19.41499755853229
$ws.Range("N21").Value = 20.34853393374712
$ws.Range("O21").Value = 32.54229300004337

$ws.Range("B22").Value = 22.60088130107779
$ws.Range("C22").Value = 8.793705205539041
$ws.Range("D22").Value = 13.32601407005705
$ws.Range("E22").Value = 13.4443852496461
$ws.Range("G22").Value = 3.699639061851406
$ws.Range("J22").Value = 8.370455386468507
$ws.Range("L22").Value = 12.28967106718243
$ws.Range("M22").Value = 19.48023643930426
$ws.Range("N22").Value = 20.29755457351417
$ws.Range("O22").Value = 32.52772711328525

$ws.Range("B23").Value = 22.45807262535684
$ws.Range("C23").Value = 8.690161296934129
$ws.Range("D23").Value = 13.32621837086752
$ws.Range("E23").Value = 13.45598534532651
$ws.Range("G23").Value = 3.700624135441525
$ws.Range("J23").Value = 8.369169798279115
$ws.Range("L23").Value = 12.28871077002685
$ws.Range("M23").Value = 19.44521648294421
$ws.Range("N23").Value = 20.32461251703695
$ws.Range("O23").Value = 32.53493681271434

$ws.Range("B24").Value = 21.91539108141151
$ws.Range("C24").Value = 8.284034918128059
$ws.Range("D24").Value = 13.33084502233461
$ws.Range("E24").Value = 13.50203894096577
$ws.Range("G24").Value = 3.704500977940826
$ws.Range("J24").Value = 8.36433170833781
$ws.Range("L24").Value = 12.28771922727599
$ws.Range("M24").Value = 19.31706538998071
$ws.Range("N24").Value = 20.43040901232304
$ws.Range("O24").Value = 32.57463226814769

$ws.Range("B25").Value = 21.32925634242483
$ws.Range("C25").Value = 7.817329274274408
$ws.Range("D25").Value = 13.34408883552019
$ws.Range("E25").Value = 13.55626371416317
$ws.Range("G25").Value = 3.708997115692263
$ws.Range("J25").Value = 8.359180519018963
$ws.Range("L25").Value = 12.29233746813402
$ws.Range("M25").Value = 19.18884744003547
$ws.Range("N25").Value = 20.55166990464361
$ws.Range("O25").Value = 32.53948653269431
